# Model training and Evaluation
# ------------------------------------------------------------------
# 1) "preprocessing" sheet: insert a new column A holding the integer
#    row index (0..4) in front of the existing label/text columns,
#    which shift from A/B to B/C.
# 2) Add a new "confusionmatrix" sheet (right after "preprocessing")
#    containing a small ham/spam confusion matrix.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- preprocessing: make room for the new index column -------------
$ws1.Columns.Item(1).Insert()

# Copy the header style (bold / bordered / centered) from the
# existing header cell onto the new index column's header-adjacent
# cells, then fill in the numeric row index.
$ws1.Range("B1").Copy()
$ws1.Range("A2:A6").PasteSpecial(-4122)

$ws1.Range("A2").Value = 0
$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2
$ws1.Range("A5").Value = 3
$ws1.Range("A6").Value = 4

# --- confusionmatrix: new sheet right after preprocessing -----------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "confusionmatrix"

$ws2.Range("B1").Value = "ham"
$ws2.Range("C1").Value = "spam"
$ws2.Range("A2").Value = "ham"
$ws2.Range("B2").Value = 1600
$ws2.Range("C2").Value = 2
$ws2.Range("A3").Value = "spam"
$ws2.Range("B3").Value = 34
$ws2.Range("C3").Value = 203

# Apply the same header style to the matrix's row/column labels.
$ws1.Range("B1").Copy()
$ws2.Range("B1:C1").PasteSpecial(-4122)
$ws2.Range("A2:A3").PasteSpecial(-4122)
